$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C10").Value = "IF Lab"
$ws.Range("C13").Value = "IF Lab"
$ws.Range("B39").Value = "IF Lab"
